# localize.xlsx update
# Adds three new localization rows to Sheet1 (key / vi / en columns):
#   - lang_unfollow   -> "Bỏ theo dõi"   / "Unfollow"
#   - lang_userName   -> "Tên tài khoản" / "Username"
#   - lang_viewMore   -> "Xem thêm…"     / "View more…"
# These land right after the existing last row (158), becoming rows 159-161,
# growing the used range from A1:C158 to A1:C161.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A159").Value = "lang_unfollow"
$ws.Range("B159").Value = "Bỏ theo dõi"
$ws.Range("C159").Value = "Unfollow"

$ws.Range("A160").Value = "lang_userName"
$ws.Range("B160").Value = "Tên tài khoản"
$ws.Range("C160").Value = "Username"

$ws.Range("A161").Value = "lang_viewMore"
$ws.Range("B161").Value = "Xem thêm…"
$ws.Range("C161").Value = "View more…"

# Match the author's final selection (active cell moved to A161).
$ws.Range("A161").Select()
